$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A13").Value = "Morado nena"
